$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Add a new "GB-ALIGNMENT" column header in G3
$ws.Cells.Item(3, 7).Value = "GB-ALIGNMENT"

# Remove the old "SingleUseId20" / "sleepy" row (row 22); this shifts
# the existing "res_txt_debug" / "Debug" row up from 23 to 22.
$ws.Rows.Item(22).Delete()

# Append the new translation rows after the existing data (now rows 23-26).
$ws.Cells.Item(23, 2).Value = "SingleUseId21"
$ws.Cells.Item(23, 3).Value = "Default"
$ws.Cells.Item(23, 4).Value = "Right"
$ws.Cells.Item(23, 5).Value = "LTR"
$ws.Cells.Item(23, 6).Value = "<result>"

$ws.Cells.Item(24, 2).Value = "SingleUseId22"
$ws.Cells.Item(24, 3).Value = "Default"
$ws.Cells.Item(24, 4).Value = "Left"
$ws.Cells.Item(24, 5).Value = "LTR"
$ws.Cells.Item(24, 6).Value = "Result"
$ws.Cells.Item(24, 7).Value = "Right"

$ws.Cells.Item(25, 2).Value = "SingleUseId23"
$ws.Cells.Item(25, 3).Value = "Default"
$ws.Cells.Item(25, 4).Value = "Left"
$ws.Cells.Item(25, 5).Value = "LTR"
$ws.Cells.Item(25, 6).Value = "<value>"

$ws.Cells.Item(26, 2).Value = "SingleUseId24"
$ws.Cells.Item(26, 3).Value = "Default"
$ws.Cells.Item(26, 4).Value = "Left"
$ws.Cells.Item(26, 5).Value = "LTR"
$ws.Cells.Item(26, 6).Value = "Debug"
